$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.032.10"
$ws.Range("E2").Value = "  -2.23%  "

$ws.Range("D3").Value = "2.503.01"
$ws.Range("E3").Value = "  -3.33%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "551.83"
$ws.Range("E5").Value = "  -3.50%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.52"
$ws.Range("E6").Value = "  -4.98%  "

$ws.Range("E7").Value = "  +0.14%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.598"
$ws.Range("E8").Value = "  -3.61%  "

$ws.Range("D9").Value = "2.502.31"
$ws.Range("E9").Value = "  -3.27%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.107"
$ws.Range("E10").Value = "  -8.98%  "

$ws.Range("E11").Value = "  -1.54%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.39"
$ws.Range("E12").Value = "  -7.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.356"
$ws.Range("E13").Value = "  -6.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.24"
$ws.Range("E14").Value = "  -6.90%  "

$ws.Range("D15").Value = "2.956.44"
$ws.Range("E15").Value = "  -3.16%  "

$ws.Range("D16").Value = "61.942.59"
$ws.Range("E16").Value = "  -1.97%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000164"
$ws.Range("E17").Value = "  -8.25%  "

$ws.Range("D18").Value = "2.505.79"
$ws.Range("E18").Value = "  -3.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.15"
$ws.Range("E19").Value = "  -6.79%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.03"
$ws.Range("E20").Value = "  -6.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.19"
$ws.Range("E21").Value = "  -7.66%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "322.92"
$ws.Range("E22").Value = "  -5.73%  "

$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.01"
$ws.Range("E24").Value = "  -4.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.75"
$ws.Range("E25").Value = "  -4.23%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000103"
$ws.Range("E26").Value = "  -4.66%  "

$ws.Range("D27").Value = "2.630.16"
$ws.Range("E27").Value = "  -3.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "544.39"
$ws.Range("E28").Value = "  -5.71%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.50"
$ws.Range("E29").Value = "  -4.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.37"
$ws.Range("E31").Value = "  -8.28%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.71"
$ws.Range("E32").Value = "  -2.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.149"
$ws.Range("E33").Value = "  -6.89%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.90"
$ws.Range("E34").Value = "  -7.41%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.58"
$ws.Range("E35").Value = "  -8.21%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.94"
$ws.Range("E36").Value = "  -9.46%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.87"
$ws.Range("E37").Value = "  -10.44%  "

$ws.Range("E38").Value = "  +0.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.380"
$ws.Range("E39").Value = "  -5.42%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.57"
$ws.Range("E40").Value = "  -5.78%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "144.01"
$ws.Range("E41").Value = "  -6.74%  "

$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.70"
$ws.Range("E43").Value = "  -8.37%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.70"
$ws.Range("E44").Value = "  -1.34%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.33"
$ws.Range("E45").Value = "  -5.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "149.61"
$ws.Range("E46").Value = "  -3.97%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.58"
$ws.Range("E47").Value = "  -8.61%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.90"
$ws.Range("E48").Value = "  -9.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0537"
$ws.Range("E49").Value = "  -8.62%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.592"
$ws.Range("E50").Value = "  -5.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0944"
$ws.Range("E51").Value = "  -6.03%  "
